$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:U39").AutoFilter(9, @("Hong Kong", "Kaohsiung"), 7)
